$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing rows whose content changes in place (style already correct) ---

# Row 16: hi_distance_full_otip
$ws.Range("A16").Value = "hi_distance_full_otip"
$ws.Range("B16").Value = 3101.856167357677
$ws.Range("C16").Value = 3350.017543582269
$ws.Range("D16").Value = -7.41
$ws.Range("E16").Value = 96.66890102150535
$ws.Range("F16").Value = 235.6447843277602
$ws.Range("G16").Value = 2949.762394878265
$ws.Range("H16").Value = 2961.264938373056
$ws.Range("I16").Value = 3206.238994396926
$ws.Range("J16").Value = 3822.019752372274

# Row 17: hsr_count_full_tip
$ws.Range("A17").Value = "hsr_count_full_tip"
$ws.Range("B17").Value = 212.2127966195529
$ws.Range("C17").Value = 197.5681176484451
$ws.Range("D17").Value = 7.41
$ws.Range("E17").Value = 11.55440508780838
$ws.Range("F17").Value = 12.53136109583805
$ws.Range("G17").Value = 195.3979334832982
$ws.Range("H17").Value = 179.8137021520933
$ws.Range("I17").Value = 223.5256608571924
$ws.Range("J17").Value = 224.3303883474233

# Row 26: highaccel_count_full_all
$ws.Range("A26").Value = "highaccel_count_full_all"
$ws.Range("B26").Value = 64.01852954116359
$ws.Range("C26").Value = 62.24959871793922
$ws.Range("D26").Value = 2.84
$ws.Range("E26").Value = 5.20065676264919
$ws.Range("F26").Value = 3.277482920213079
$ws.Range("G26").Value = 55.23046715121213
$ws.Range("H26").Value = 56.6623185117708
$ws.Range("I26").Value = 68.60283618067439
$ws.Range("J26").Value = 66.92139647224181

# Row 27: total_metersperminute_full_otip
$ws.Range("A27").Value = "total_metersperminute_full_otip"
$ws.Range("B27").Value = 1854.429759905351
$ws.Range("C27").Value = 1898.640765806154
$ws.Range("D27").Value = -2.33
$ws.Range("E27").Value = 80.93456021654049
$ws.Range("F27").Value = 53.94696992923298
$ws.Range("G27").Value = 1710.777219034276
$ws.Range("H27").Value = 1771.447233404109
$ws.Range("I27").Value = 1900.827443428021
$ws.Range("J27").Value = 1990.21486077696

# --- New rows appended below the previous last row (need style copied onto column A) ---
$styleSrc = $ws.Range("A2")

# Row 28: highdecel_count_full_all
$styleSrc.Copy($ws.Range("A28"))
$ws.Range("A28").Value = "highdecel_count_full_all"
$ws.Range("B28").Value = 150.4994073400977
$ws.Range("C28").Value = 153.4089909332751
$ws.Range("D28").Value = -1.9
$ws.Range("E28").Value = 9.65149590027513
$ws.Range("F28").Value = 6.048299869578941
$ws.Range("G28").Value = 138.3932627864787
$ws.Range("H28").Value = 145.4461546644024
$ws.Range("I28").Value = 162.0473069603744
$ws.Range("J28").Value = 169.212898593574

# Row 29: medaccel_count_full_all
$styleSrc.Copy($ws.Range("A29"))
$ws.Range("A29").Value = "medaccel_count_full_all"
$ws.Range("B29").Value = 1036.710149652919
$ws.Range("C29").Value = 1054.171502221885
$ws.Range("D29").Value = -1.66
$ws.Range("E29").Value = 13.13270641296326
$ws.Range("F29").Value = 18.35008164869249
$ws.Range("G29").Value = 1026.938123561755
$ws.Range("H29").Value = 1020.74291875804
$ws.Range("I29").Value = 1051.554767565639
$ws.Range("J29").Value = 1088.759401937111

# Row 30: sprint_count_full_all
$styleSrc.Copy($ws.Range("A30"))
$ws.Range("A30").Value = "sprint_count_full_all"
$ws.Range("B30").Value = 88.43482101811608
$ws.Range("C30").Value = 89.84853943724046
$ws.Range("D30").Value = -1.57
$ws.Range("E30").Value = 6.504426500192894
$ws.Range("F30").Value = 5.545015069176541
$ws.Range("G30").Value = 81.59013494508247
$ws.Range("H30").Value = 79.46233176170961
$ws.Range("I30").Value = 95.18758333784231
$ws.Range("J30").Value = 99.89499369977753

# Row 31: total_metersperminute_full_tip
$styleSrc.Copy($ws.Range("A31"))
$ws.Range("A31").Value = "total_metersperminute_full_tip"
$ws.Range("B31").Value = 1759.309419383992
$ws.Range("C31").Value = 1785.157316325987
$ws.Range("D31").Value = -1.45
$ws.Range("E31").Value = 68.34387076500634
$ws.Range("F31").Value = 59.49952442158469
$ws.Range("G31").Value = 1642.562416926959
$ws.Range("H31").Value = 1697.513342961181
$ws.Range("I31").Value = 1813.472520353595
$ws.Range("J31").Value = 1917.084750358108

# Row 32: psv99
$styleSrc.Copy($ws.Range("A32"))
$ws.Range("A32").Value = "psv99"
$ws.Range("B32").Value = 358.3587389697205
$ws.Range("C32").Value = 362.6029586435135
$ws.Range("D32").Value = -1.17
$ws.Range("E32").Value = 16.1618311469886
$ws.Range("F32").Value = 8.889284818201046
$ws.Range("G32").Value = 337.581470170531
$ws.Range("H32").Value = 349.2217523827069
$ws.Range("I32").Value = 376.9914152822696
$ws.Range("J32").Value = 379.6420174970298

# Row 33: total_metersperminute_full_all
$styleSrc.Copy($ws.Range("A33"))
$ws.Range("A33").Value = "total_metersperminute_full_all"
$ws.Range("B33").Value = 1467.915128687428
$ws.Range("C33").Value = 1482.665020609128
$ws.Range("D33").Value = -0.99
$ws.Range("E33").Value = 46.88355439183032
$ws.Range("F33").Value = 41.14863069641145
$ws.Range("G33").Value = 1398.95043196687
$ws.Range("H33").Value = 1414.563064016289
$ws.Range("I33").Value = 1520.66364081299
$ws.Range("J33").Value = 1558.253250006369

# Row 34: meddecel_count_full_all
$styleSrc.Copy($ws.Range("A34"))
$ws.Range("A34").Value = "meddecel_count_full_all"
$ws.Range("B34").Value = 856.3854211848644
$ws.Range("C34").Value = 863.9768443859568
$ws.Range("D34").Value = -0.88
$ws.Range("E34").Value = 13.84175151793503
$ws.Range("F34").Value = 13.45453799333878
$ws.Range("G34").Value = 839.751211908086
$ws.Range("H34").Value = 831.952186859352
$ws.Range("I34").Value = 870.5006003070322
$ws.Range("J34").Value = 885.0464636511241

# Row 35: hsr_distance_full_all
$styleSrc.Copy($ws.Range("A35"))
$ws.Range("A35").Value = "hsr_distance_full_all"
$ws.Range("B35").Value = 5152.708273705066
$ws.Range("C35").Value = 5194.68376096633
$ws.Range("D35").Value = -0.8100000000000001
$ws.Range("E35").Value = 197.4592798991305
$ws.Range("F35").Value = 212.2926671275888
$ws.Range("G35").Value = 4946.909167362483
$ws.Range("H35").Value = 4887.994623704785
$ws.Range("I35").Value = 5375.009235331764
$ws.Range("J35").Value = 5603.880245002479

# Row 36: hi_distance_full_all
$styleSrc.Copy($ws.Range("A36"))
$ws.Range("A36").Value = "hi_distance_full_all"
$ws.Range("B36").Value = 6684.369914484892
$ws.Range("C36").Value = 6733.326361891136
$ws.Range("D36").Value = -0.73
$ws.Range("E36").Value = 345.6754561857296
$ws.Range("F36").Value = 281.5000371572231
$ws.Range("G36").Value = 6323.785658895358
$ws.Range("H36").Value = 6316.044879787938
$ws.Range("I36").Value = 7043.33210958169
$ws.Range("J36").Value = 7175.391784482549

# Row 37: hi_count_full_all
$styleSrc.Copy($ws.Range("A37"))
$ws.Range("A37").Value = "hi_count_full_all"
$ws.Range("B37").Value = 586.4490446983772
$ws.Range("C37").Value = 589.789154483561
$ws.Range("D37").Value = -0.57
$ws.Range("E37").Value = 22.06536729230082
$ws.Range("F37").Value = 22.31730342174017
$ws.Range("G37").Value = 557.5263776349524
$ws.Range("H37").Value = 558.7246739002793
$ws.Range("I37").Value = 612.3824791739308
$ws.Range("J37").Value = 630.2000857110838

# Row 38: sprint_distance_full_all
$styleSrc.Copy($ws.Range("A38"))
$ws.Range("A38").Value = "sprint_distance_full_all"
$ws.Range("B38").Value = 1531.661640779826
$ws.Range("C38").Value = 1538.642600924807
$ws.Range("D38").Value = -0.45
$ws.Range("E38").Value = 175.2120404592917
$ws.Range("F38").Value = 100.6629281253413
$ws.Range("G38").Value = 1311.192578310617
$ws.Range("H38").Value = 1333.496865591414
$ws.Range("I38").Value = 1691.404988595376
$ws.Range("J38").Value = 1697.329065606786

# Row 39: hsr_count_full_all
$styleSrc.Copy($ws.Range("A39"))
$ws.Range("A39").Value = "hsr_count_full_all"
$ws.Range("B39").Value = 498.0142236802611
$ws.Range("C39").Value = 499.9406150463204
$ws.Range("D39").Value = -0.39
$ws.Range("E39").Value = 15.88777695840753
$ws.Range("F39").Value = 17.12807449190218
$ws.Range("G39").Value = 475.9362426898699
$ws.Range("H39").Value = 479.2623421385697
$ws.Range("I39").Value = 517.1948958360886
$ws.Range("J39").Value = 530.3050920113063

# Row 40: running_distance_full_all
$styleSrc.Copy($ws.Range("A40"))
$ws.Range("A40").Value = "running_distance_full_all"
$ws.Range("B40").Value = 13841.85545268292
$ws.Range("C40").Value = 13893.78291834627
$ws.Range("D40").Value = -0.37
$ws.Range("E40").Value = 631.0648817821947
$ws.Range("F40").Value = 580.2522196518806
$ws.Range("G40").Value = 13119.47707968072
$ws.Range("H40").Value = 13116.43486990435
$ws.Range("I40").Value = 14440.29166203201
$ws.Range("J40").Value = 14949.86472154682

# Row 41: total_distance_full_all
$styleSrc.Copy($ws.Range("A41"))
$ws.Range("A41").Value = "total_distance_full_all"
$ws.Range("B41").Value = 101600.5653035157
$ws.Range("C41").Value = 101949.6126464969
$ws.Range("D41").Value = -0.34
$ws.Range("E41").Value = 639.822932365716
$ws.Range("F41").Value = 1631.687361247225
$ws.Range("G41").Value = 101121.5096110575
$ws.Range("H41").Value = 99233.95272570745
$ws.Range("I41").Value = 102631.618022524
$ws.Range("J41").Value = 104631.0950484905

Write-Output "done"